# Adds the "Valid Parenthesis" problem as a new row to the Leetcode Notes
# table, and adds a new "Notes" column capturing a note about the solution.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Grow the table by one column (Table1 goes from A1:E162 to A1:F162).
# Excel auto-names the new column "ColumnN"; it gets its real name below
# once we write the header text into F1.
$null = $tbl.ListColumns.Add()

# Row 5 already has "#" (A5 = 4) filled in; fill in the rest of the new
# "Valid Parenthesis" entry.
$ws.Range("B5").Value = "Valid Parenthesis"
$ws.Range("C5").Value = "Given a string, return whether paranthesis are valid (meets 3 conditions)"
$ws.Range("D5").Value = "Using a stack, pushing when open paranthesis are found and popping with close paranthesis"

# Header for the newly added table column.
$ws.Range("F1").Value = "Notes"

$ws.Range("E5").Value = 8
$ws.Range("F5").Value = "Logic itself was easy, mostly just syntax issues"

# Match Excel's natural cursor position after finishing data entry on row 5
# (Enter moves the active cell down into the next row).
$null = $ws.Range("B6").Select()
